$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "time_taken" in column F, row 1 (mirrors style/format of existing header cells)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Map of row number -> time_taken value (as captured from the source data)
$timeTaken = @{
    2  = "2021-10-05 10:51:30.064848"
    3  = "2021-10-05 10:51:30.064861"
    4  = "2021-10-05 10:51:30.064865"
    5  = "2021-10-05 10:51:30.064869"
    6  = "2021-10-05 10:51:30.064872"
    7  = "2021-10-05 10:51:30.064875"
    8  = "2021-10-05 10:51:30.064878"
    9  = "2021-10-05 10:51:30.064881"
    10 = "2021-10-05 10:51:30.064885"
    11 = "2021-10-05 10:51:30.064888"
    12 = "2021-10-05 10:51:30.064891"
    13 = "2021-10-05 10:51:30.064894"
    14 = "2021-10-05 10:51:30.064897"
    15 = "2021-10-05 10:51:30.064901"
    16 = "2021-10-05 10:51:30.064904"
    17 = "2021-10-05 10:51:30.064907"
    18 = "2021-10-05 10:51:30.064910"
    19 = "2021-10-05 10:51:30.064913"
    20 = "2021-10-05 10:51:30.064916"
    21 = "2021-10-05 10:51:30.064919"
    22 = "2021-10-05 10:51:30.064922"
    23 = "2021-10-05 10:51:30.064926"
    24 = "2021-10-05 10:51:30.064929"
    25 = "2021-10-05 10:51:30.064932"
    26 = "2021-10-05 10:51:30.064935"
    27 = "2021-10-05 10:51:30.064938"
    28 = "2021-10-05 10:51:30.064941"
    29 = "2021-10-05 10:51:30.064944"
    30 = "2021-10-05 10:51:30.064947"
}

foreach ($row in $timeTaken.Keys) {
    $ws.Cells.Item($row, 6).Value = $timeTaken[$row]
}
